# Update metrics_6_6 sheet: refresh all metric columns (B:Q) with the new
# ensemble-wide values, and reorder the model names in column A according
# to the new training run ("atualizado todo o treinamento para o novo lm").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New model name order for rows 2..26 (row -> model name)
$newNames = @{
    2  = "model_6_6_0"
    3  = "model_6_6_22"
    4  = "model_6_6_21"
    5  = "model_6_6_20"
    6  = "model_6_6_19"
    7  = "model_6_6_18"
    8  = "model_6_6_17"
    9  = "model_6_6_16"
    10 = "model_6_6_15"
    11 = "model_6_6_14"
    12 = "model_6_6_13"
    13 = "model_6_6_23"
    14 = "model_6_6_12"
    15 = "model_6_6_10"
    16 = "model_6_6_9"
    17 = "model_6_6_8"
    18 = "model_6_6_7"
    19 = "model_6_6_6"
    20 = "model_6_6_5"
    21 = "model_6_6_4"
    22 = "model_6_6_3"
    23 = "model_6_6_2"
    24 = "model_6_6_1"
    25 = "model_6_6_11"
    26 = "model_6_6_24"
}

# New shared metric values (B:Q) applied to every model row
# (written as plain decimals; this PowerShell parser does not accept
# scientific-notation literals like 3.43e-05)
$newValues = @{
    "B" = 0.9999632043074972
    "C" = 0.9992467334840777
    "D" = 0.999963735509305
    "E" = 0.9999954468363915
    "F" = 0.9999748712582449
    "G" = 0.00003434716380714357
    "H" = 0.0007031412280348575
    "I" = 0.0000174263632936475
    "J" = 0.0000003328705626383837
    "K" = 0.000008879616928142945
    "L" = 0.0002886930549261025
    "M" = 0.005860645340501639
    "N" = 1.000883096620067
    "O" = 0.006110145145855123
    "P" = 70.55798222036485
    "Q" = 101.0298778420699
}

foreach ($row in 2..26) {
    $ws.Range("A$row").Value = $newNames[$row]
    foreach ($col in $newValues.Keys) {
        $ws.Range("$col$row").Value = $newValues[$col]
    }
}
